# Insert three new rows at 287 (pushing the existing 287-371 block down to
# 290-374), then populate the three new rows with the new data point
# (Hayward Kiwi, 2021-09-29, Provincia de Curico, $/bandeja 10 kilos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("287:289").Insert()

# --- Row 287: Especial ---
$ws.Cells.Item(287,1).Value2 = 9
$ws.Cells.Item(287,2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(287,3).Value2 = "Metropolitana"
$ws.Cells.Item(287,4).Value2 = 44468
$ws.Cells.Item(287,5).Value2 = 13
$ws.Cells.Item(287,6).Value2 = "Fruta"
$ws.Cells.Item(287,7).Value2 = 100101
$ws.Cells.Item(287,8).Value2 = "Berries"
$ws.Cells.Item(287,9).Value2 = 100101007
$ws.Cells.Item(287,10).Value2 = "Kiwi"
$ws.Cells.Item(287,11).Value2 = "Hayward"
$ws.Cells.Item(287,12).Value2 = "Especial"
$ws.Cells.Item(287,13).Value2 = 300
$ws.Cells.Item(287,14).Value2 = 11000
$ws.Cells.Item(287,15).Value2 = 11000
$ws.Cells.Item(287,16).Value2 = 11000
$ws.Cells.Item(287,17).Value2 = "`$/bandeja 10 kilos"
$ws.Cells.Item(287,18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(287,19).Value2 = 1100
$ws.Cells.Item(287,20).Value2 = 10

# --- Row 288: Primera ---
$ws.Cells.Item(288,1).Value2 = 9
$ws.Cells.Item(288,2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(288,3).Value2 = "Metropolitana"
$ws.Cells.Item(288,4).Value2 = 44468
$ws.Cells.Item(288,5).Value2 = 13
$ws.Cells.Item(288,6).Value2 = "Fruta"
$ws.Cells.Item(288,7).Value2 = 100101
$ws.Cells.Item(288,8).Value2 = "Berries"
$ws.Cells.Item(288,9).Value2 = 100101007
$ws.Cells.Item(288,10).Value2 = "Kiwi"
$ws.Cells.Item(288,11).Value2 = "Hayward"
$ws.Cells.Item(288,12).Value2 = "Primera"
$ws.Cells.Item(288,13).Value2 = 300
$ws.Cells.Item(288,14).Value2 = 9000
$ws.Cells.Item(288,15).Value2 = 9000
$ws.Cells.Item(288,16).Value2 = 9000
$ws.Cells.Item(288,17).Value2 = "`$/bandeja 10 kilos"
$ws.Cells.Item(288,18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(288,19).Value2 = 900
$ws.Cells.Item(288,20).Value2 = 10

# --- Row 289: Segunda ---
$ws.Cells.Item(289,1).Value2 = 9
$ws.Cells.Item(289,2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(289,3).Value2 = "Metropolitana"
$ws.Cells.Item(289,4).Value2 = 44468
$ws.Cells.Item(289,5).Value2 = 13
$ws.Cells.Item(289,6).Value2 = "Fruta"
$ws.Cells.Item(289,7).Value2 = 100101
$ws.Cells.Item(289,8).Value2 = "Berries"
$ws.Cells.Item(289,9).Value2 = 100101007
$ws.Cells.Item(289,10).Value2 = "Kiwi"
$ws.Cells.Item(289,11).Value2 = "Hayward"
$ws.Cells.Item(289,12).Value2 = "Segunda"
$ws.Cells.Item(289,13).Value2 = 250
$ws.Cells.Item(289,14).Value2 = 7000
$ws.Cells.Item(289,15).Value2 = 7000
$ws.Cells.Item(289,16).Value2 = 7000
$ws.Cells.Item(289,17).Value2 = "`$/bandeja 10 kilos"
$ws.Cells.Item(289,18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(289,19).Value2 = 700
$ws.Cells.Item(289,20).Value2 = 10
